# Append: 2025-12-03 18:37 JST
# Update the "取得日時" (acquisition datetime) column (A) for all existing
# data rows (rows 2-13) on the "ランサーズ" sheet from the old timestamp
# "2025-12-03 18:29:31" to the new timestamp "2025-12-03 18:37:34".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldTimestamp = "2025-12-03 18:29:31"
$newTimestamp = "2025-12-03 18:37:34"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}
